$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Hunk 1: "Information will not [_GoBack bookmark] display at URL and
# using lag." -> remove the _GoBack bookmark and merge the two runs
# into a single run of plain text.
# -----------------------------------------------------------------

$null = $d.Content.Find.Execute(
    "Information will not display at URL and using lag.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Information will not display at URL and using lag.", 2)

# -----------------------------------------------------------------
# Hunk 2: "it’s" (wrapped in spellStart/gramStart .. spellEnd/gramEnd
# proofErr markers) -> "its" (proofErr markers removed), and a new
# _GoBack bookmark is inserted right after "its".
# -----------------------------------------------------------------

# Locate the word "it's" (curly apostrophe) together with one
# character of context on each side, so that replacing the span
# merges it into its neighbouring runs and drops the now orphaned
# <w:proofErr/> markers that surrounded it.
$findCtx = $d.Content.Duplicate
$found = $findCtx.Find.Execute("d it" + [char]8217 + "s c")
if (-not $found) {
    throw "Could not find the 'it's' context span"
}
$ctxStart = $findCtx.Start
$ctxEnd = $findCtx.End

$wideRange = $d.Range($ctxStart, $ctxEnd)
$wideRange.Text = "d its c"

# Re-locate "and its" so we can split "and " away from "its" into its
# own run (matching the target run layout), and then again so that we
# can insert the _GoBack bookmark right after "its".
$andIts = $d.Content.Duplicate
$null = $andIts.Find.Execute("resource and")
$splitPos1 = $andIts.End + 1
$d.Bookmarks.Add("TempSplit", $d.Range($splitPos1, $splitPos1)) | Out-Null
$d.Bookmarks("TempSplit").Delete()

$itsRange = $d.Content.Duplicate
$null = $itsRange.Find.Execute("and its")
$bmPos = $itsRange.End
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos)) | Out-Null
